# Add a "Current Fiscal Year" row (row 10) to Sheet1, between the existing
# "Contact Info" block (rows 7-9) and the "Overall Numbers" section (row 11),
# mirroring the style used by the other question rows in that block (A7:A9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the formatting of the row above (A9) down into the new row (A10) so the
# new label picks up the same style used by its neighboring question rows.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)

$ws.Range("A10").Value = "Current Fiscal Year"
$ws.Range("B10").Value = 2023
